# Update the CDA Logical Model metadata sheet (BL: Boolean) for the
# ST.r2b regeneration:
#   - bump the Version property
#   - bump the Date property
#   - insert a new "Jurisdiction" property row right after "Contact"
#
# The workbook has two worksheets:
#   1 = "Metadata" (Property/Value pairs, one per row)
#   2 = "Elements" (the structure definition element table)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Update the "Version" value (row 3, column B) -----------------------
$ws.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"

# --- 2. Update the "Date" value (row 8, column B) ---------------------------
$ws.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# --- 3. Insert a new "Jurisdiction" / "" row right after "Contact" ---------
# "Contact" currently lives on row 10, and rows 11..19 hold the remaining
# properties (Description .. Derivation). Shift them all down one row using
# Copy (keeps the original cell typing, e.g. text "false" stays text instead
# of turning into a boolean), clearing the destination first so that blank
# source cells actually blank out the destination instead of leaving stale
# values behind.
for ($r = 19; $r -ge 11; $r--) {
    $src = $ws.Range("A" + $r + ":B" + $r)
    $dstTopLeft = $ws.Range("A" + ($r + 1))
    $dstFull = $ws.Range("A" + ($r + 1) + ":B" + ($r + 1))
    $dstFull.ClearContents()
    $src.Copy($dstTopLeft)
}

# Row 11 is now free for the new "Jurisdiction" property (value left blank).
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
